$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = $origStyle
}

Set-TextValue "D2" "310.26"
Set-TextValue "E2" "0.53%"
Set-TextValue "D3" "41.02"
Set-TextValue "E3" "-0.37%"
Set-TextValue "D4" "5.217"
Set-TextValue "E4" "1.57%"
Set-TextValue "D5" "0.07682"
Set-TextValue "E5" "0.57%"
Set-TextValue "D6" "4.281"
Set-TextValue "E6" "0.26%"
Set-TextValue "D7" "1.738"
Set-TextValue "E7" "7.59%"
Set-TextValue "D8" "0.9303"
Set-TextValue "E8" "2.64%"
Set-TextValue "D10" "0.1281"
Set-TextValue "E10" "13.08%"
Set-TextValue "D11" "0.1830"
Set-TextValue "E11" "1.87%"
Set-TextValue "D12" "0.09126"
Set-TextValue "E12" "0.09%"
Set-TextValue "D13" "0.04214"
Set-TextValue "E13" "-0.42%"
Set-TextValue "D14" "0.1052"
Set-TextValue "E14" "0.23%"
Set-TextValue "D15" "0.001290"
Set-TextValue "E15" "2.45%"
Set-TextValue "D16" "0.005887"
Set-TextValue "E16" "2.76%"
Set-TextValue "E17" "0.21%"
Set-TextValue "D19" "7.389"
Set-TextValue "E19" "9.84%"
Set-TextValue "D20" "0.1351"
Set-TextValue "E20" "-1.03%"
Set-TextValue "D21" "0.2717"
Set-TextValue "E21" "-0.75%"
Set-TextValue "D22" "0.04018"
Set-TextValue "E22" "-1.42%"
Set-TextValue "E23" "-0.03%"
Set-TextValue "D24" "0.004099"
Set-TextValue "E24" "1.43%"
Set-TextValue "D25" "0.0001271"
Set-TextValue "E25" "0.02%"
Set-TextValue "D38" "0.02536"
Set-TextValue "E38" "4.71%"
Set-TextValue "D39" "0.05322"
Set-TextValue "E39" "1.62%"
Set-TextValue "D40" "0.007857"
Set-TextValue "E40" "0.64%"
Set-TextValue "D41" "0.1313"
Set-TextValue "E41" "0.89%"
Set-TextValue "D42" "0.006649"
Set-TextValue "E42" "1.71%"
Set-TextValue "D43" "0.002053"
Set-TextValue "E43" "5.25%"
Set-TextValue "D44" "0.008093"
Set-TextValue "E44" "6.90%"
Set-TextValue "D45" "0.3084"
Set-TextValue "E45" "-0.10%"
Set-TextValue "D46" "0.00006789"
Set-TextValue "E46" "0.01%"
Set-TextValue "E47" "0.02%"
Set-TextValue "D48" "0.2248"
Set-TextValue "E48" "206.30%"
Set-TextValue "D50" "0.00002101"
Set-TextValue "E50" "0.02%"
Set-TextValue "D51" "0.0002001"
Set-TextValue "E51" "0.02%"
